# Update cryptos list (prices + 1h volume change) as scraped on
# Fri Apr 19 00:00:30 UTC 2024 with GitHub Actions.
# Values that look like plain numbers are prefixed with a leading
# apostrophe so Excel keeps them as text (matching the source data,
# which preserves trailing zeros such as "7.90" or "6.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.459.61"
$ws.Range("E2").Value = "  +3.51%  "
$ws.Range("D3").Value = "3.064.05"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'551.72"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "'141.98"
$ws.Range("E6").Value = "  +7.33%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.060.50"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  +6.82%  "
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").Value = "'34.81"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "3.563.52"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "63.437.72"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "3.066.42"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").Value = "'482.78"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'13.88"
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  +5.86%  "
$ws.Range("D24").Value = "'80.76"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("E25").Value = "  +7.55%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +7.26%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("E33").Value = "  +7.97%  "
$ws.Range("E34").Value = "  +5.09%  "
$ws.Range("D35").Value = "'55.28"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "'465.06"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("D39").Value = "'0.0396"
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.119"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.999.52"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("E43").Value = "  +5.67%  "
$ws.Range("E44").Value = "  +5.90%  "
$ws.Range("D45").Value = "'0.256"
$ws.Range("E45").Value = "  +5.69%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'2.05"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "'116.47"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("E51").Value = "  +4.44%  "
